# Update "want-to-go" counts (column F) for a handful of events that
# appear on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 49
$ws1.Range("F4").Value  = 619
$ws1.Range("F5").Value  = 159
$ws1.Range("F6").Value  = 9345
$ws1.Range("F7").Value  = 840
$ws1.Range("F9").Value  = 1191
$ws1.Range("F10").Value = 1116
$ws1.Range("F11").Value = 145
$ws1.Range("F12").Value = 89
$ws1.Range("F13").Value = 16
$ws1.Range("F14").Value = 257
$ws1.Range("F15").Value = 405
$ws1.Range("F18").Value = 1251

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 49
$ws4.Range("F5").Value  = 619
$ws4.Range("F6").Value  = 159
$ws4.Range("F7").Value  = 9345
$ws4.Range("F8").Value  = 840
$ws4.Range("F10").Value = 1191
$ws4.Range("F11").Value = 1116
$ws4.Range("F12").Value = 145
$ws4.Range("F13").Value = 89
$ws4.Range("F14").Value = 16
$ws4.Range("F15").Value = 257
$ws4.Range("F16").Value = 405
$ws4.Range("F19").Value = 1251
